$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for new columns I and J
# Copy formatting from existing header cell (H1) so borders/bold/alignment match,
# then overwrite with the new header text.
$ws.Range("H1").Copy($ws.Range("I1"))
$ws.Range("H1").Copy($ws.Range("J1"))
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for columns I (I0) and J (IF), rows 2-22
$values = @(
    @(10, 10),
    @(7, 7),
    @(8, 8),
    @(7, 7),
    @(3, 4),
    @(7, 7),
    @(7, 7),
    @(4, 4),
    @(9, 9),
    @(8, 9),
    @(9, 9),
    @(2, 4),
    @(7, 8),
    @(8, 8),
    @(8, 8),
    @(7, 8),
    @(7, 8),
    @(5, 6),
    @(7, 8),
    @(8, 8),
    @(8, 8)
)

$row = 2
foreach ($pair in $values) {
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
    $row++
}
